$d = $word.ActiveDocument

# Locate the start of the Heading3 paragraph "Trener unosi lozinku koja
# nije jedinstvena u sistemu" and delete it together with the two
# following list paragraphs that describe that (now removed) scenario:
#   - "Koraci 1 i 2 su isti kao u slucaju 2.2.1."
#   - " Trener pritiska dugme "Commit Changes". Treneru iskace poruka
#      o tome da uneta lozinka vec postoji u bazi, ne dolazi do promene
#      u sistemu. "
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($startPara -eq $null -and $t -like "*Trener unosi lozinku koja nije jedinstvena*") {
        $startPara = $p
    }
    elseif ($startPara -ne $null -and $endPara -eq $null -and $t -like "*Trener pritiska dugme*") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
